$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-11 Wednesday" "2024-09-12 Thursday"

Replace-Text "796÷2=398, 0" "254÷3=84, 2"
Replace-Text "150÷5=30, 0" "778÷7=111, 1"
Replace-Text "493÷5=98, 3" "809÷3=269, 2"
Replace-Text "699÷4=174, 3" "566÷2=283, 0"
Replace-Text "737÷5=147, 2" "828÷5=165, 3"

Replace-Text "245÷7=35, 0" "564÷9=62, 6"
Replace-Text "531÷7=75, 6" "408÷2=204, 0"
Replace-Text "290÷5=58, 0" "138÷6=23, 0"
Replace-Text "172÷3=57, 1" "694÷9=77, 1"
Replace-Text "983÷8=122, 7" "332÷8=41, 4"

Replace-Text "418÷8=52, 2" "561÷2=280, 1"
Replace-Text "695÷5=139, 0" "662÷9=73, 5"
Replace-Text "672÷6=112, 0" "114÷7=16, 2"
Replace-Text "573÷6=95, 3" "733÷6=122, 1"
Replace-Text "506÷8=63, 2" "307÷4=76, 3"

Replace-Text "451÷9=50, 1" "882÷2=441, 0"
Replace-Text "569÷5=113, 4" "401÷8=50, 1"
Replace-Text "297÷8=37, 1" "215÷6=35, 5"
Replace-Text "402÷2=201, 0" "117÷9=13, 0"
Replace-Text "728÷5=145, 3" "147÷3=49, 0"

Replace-Text "301÷2=150, 1" "281÷4=70, 1"
Replace-Text "307÷9=34, 1" "955÷3=318, 1"
Replace-Text "265÷8=33, 1" "902÷2=451, 0"
Replace-Text "426÷7=60, 6" "587÷7=83, 6"
Replace-Text "104÷6=17, 2" "726÷2=363, 0"

Write-Host "Done"
